$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: clone the visible formatting of a "template" cell onto a target
# cell (font, fill, number format). The workbook's existing rows follow a
# repeating 2-row pattern (a "field name" row followed by a "value" row), so
# the four new rows being appended (8-11) reuse the same look as the most
# recent existing pair (rows 6-7).
# ---------------------------------------------------------------------------
function Copy-CellLook($srcAddr, $dstAddr) {
    $src = $ws.Range($srcAddr)
    $dst = $ws.Range($dstAddr)
    $dst.Font.Name = $src.Font.Name
    $dst.Font.Size = $src.Font.Size
    $dst.Font.Bold = $src.Font.Bold
    $dst.Font.Color = $src.Font.Color
    $dst.Interior.Pattern = $src.Interior.Pattern
    if ($src.Interior.Pattern -ne -4142) {
        $dst.Interior.Color = $src.Interior.Color
    }
    $dst.NumberFormat = $src.NumberFormat
}

# ---------------------------------------------------------------------------
# Row 8: "field name" row (same look as rows 2/4/6)
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "SANITY5"
$ws.Range("B8").Value = "ExploreSearch"
$ws.Range("C8").Value = "ProcessName"
$ws.Range("D8").Value = "ExploreChildSearch"
$ws.Range("E8").Value = "IssueText"
$ws.Range("F8").Value = "BillingType"

Copy-CellLook "A6" "A8"
Copy-CellLook "B6" "B8"
Copy-CellLook "C6" "C8"
Copy-CellLook "D6" "D8"
Copy-CellLook "E6" "E8"
Copy-CellLook "F6" "F8"
$ws.Rows(8).RowHeight = 15

# ---------------------------------------------------------------------------
# Row 9: "value" row (same look as row 7)
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = "DATA_SANITY5"
$ws.Range("B9").Value = "AUTOMATION SEARCH"
$ws.Range("C9").Value = "scenario5"
$ws.Range("D9").Value = "Work Orders"
$ws.Range("E9").Value = "Billing Type cannot be Loan"
$ws.Range("F9").Value = "Loan"

Copy-CellLook "B7" "A9"
Copy-CellLook "B7" "B9"
Copy-CellLook "C7" "C9"
Copy-CellLook "D7" "D9"
Copy-CellLook "E7" "E9"
Copy-CellLook "E7" "F9"
$ws.Rows(9).RowHeight = 15

# ---------------------------------------------------------------------------
# Row 10: "field name" row (same look as row 8 above, but with only the
# first four columns carrying text and the last two left blank)
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "SANITY6"
$ws.Range("B10").Value = "ExploreSearch"
$ws.Range("C10").Value = "ProcessName"
$ws.Range("D10").Value = "ExploreChildSearch"

Copy-CellLook "A8" "A10"
Copy-CellLook "B8" "B10"
Copy-CellLook "C8" "C10"
Copy-CellLook "D8" "D10"
Copy-CellLook "E8" "E10"
Copy-CellLook "F8" "F10"
$ws.Rows(10).RowHeight = 15

# ---------------------------------------------------------------------------
# Row 11: "value" row (same look as row 9 above), last two columns blank
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = "DATA_SANITY6"
$ws.Range("B11").Value = "AUTOMATION SEARCH"
$ws.Range("C11").Value = "scenario5"
$ws.Range("D11").Value = "Cases"

Copy-CellLook "A9" "A11"
Copy-CellLook "B9" "B11"
Copy-CellLook "C9" "C11"
Copy-CellLook "D9" "D11"
Copy-CellLook "E9" "E11"
Copy-CellLook "F9" "F11"
$ws.Rows(11).RowHeight = 15

# ---------------------------------------------------------------------------
# A7 is re-aligned to the same look as the rest of its row (B7:E7)
# ---------------------------------------------------------------------------
Copy-CellLook "B7" "A7"

# ---------------------------------------------------------------------------
# Column E widens slightly to fit the new "IssueText"/issue-description data
# ---------------------------------------------------------------------------
$ws.Columns("E:E").ColumnWidth = 21.45

# ---------------------------------------------------------------------------
# Update the active selection to reflect where editing left off
# ---------------------------------------------------------------------------
$ws.Range("C15").Select() | Out-Null
